$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Add the new BOM row (Schakelaar / Aan-uit tuimelschakelaar)
$ws.Range("B18").Value = "Schakelaar"
$ws.Range("C18").Value = "Aan/uit tuimelschakelaar"

# Match the author's final selection/scroll position
$ws.Range("C18").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
